$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("B3").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("B7").Value = 0
$ws.Range("B15").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("B17").Value = 50
$ws.Range("B21").Value = 0
$ws.Range("B23").Value = 5
$ws.Range("B24").Value = 5
$ws.Range("D26").Value = 5

$ws.Range("B18").Select()
